$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18; this pushes the existing rows 18-85
# down to 19-86 (matching the row-shift pattern seen across the rest of
# the sheet).
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new data record.
$ws.Cells.Item(18, 1).Value = 11
$ws.Cells.Item(18, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(18, 3).Value = "Bíobío"
$ws.Cells.Item(18, 4).Value = 44624
$ws.Cells.Item(18, 5).Value = 8
$ws.Cells.Item(18, 6).Value = 100112024
$ws.Cells.Item(18, 7).Value = "Choclo"
$ws.Cells.Item(18, 8).Value = "Choclero"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 27000
$ws.Cells.Item(18, 11).Value = 150
$ws.Cells.Item(18, 12).Value = 170
$ws.Cells.Item(18, 13).Value = 159
$ws.Cells.Item(18, 14).Value = "$/unidad"
$ws.Cells.Item(18, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(18, 16).Value = 159
$ws.Cells.Item(18, 17).Value = 1
$ws.Cells.Item(18, 18).Value = "Hortaliza"
